$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 1349.0857
$ws.Range("I17").Value = 799
$ws.Range("J17").Value = 1365.2646
$ws.Range("K17").Value = 2397
$ws.Range("L17").Value = 4095.7938
$ws.Range("M17").Value = -2229
$ws.Range("N17").Value = -4431.793799999999

# row 19
$ws.Range("H19").Value = 1214.8462
$ws.Range("I19").Value = 382
$ws.Range("J19").Value = 1585
$ws.Range("K19").Value = 382
$ws.Range("L19").Value = 1585
$ws.Range("M19").Value = -207
$ws.Range("N19").Value = -1935

# row 33
$ws.Range("H33").Value = 80
$ws.Range("I33").Value = 80
$ws.Range("K33").Value = 80
$ws.Range("M33").Value = 149

# row 40
$ws.Range("H40").Value = 11112844
$ws.Range("I40").Value = 2600
$ws.Range("J40").Value = 33333332
$ws.Range("K40").Value = 2600
$ws.Range("L40").Value = 33333332
$ws.Range("M40").Value = -2425
$ws.Range("N40").Value = -33333682

# row 58
$ws.Range("H58").Value = 621.625
$ws.Range("J58").Value = 1999.5
$ws.Range("L58").Value = 5998.5
$ws.Range("N58").Value = -6298.5

# row 62
$ws.Range("H62").Value = 4006
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4006
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4006
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5254

# row 65
$ws.Range("H65").Value = 4006
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4006
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20030
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26270

# row 75
$ws.Range("H75").Value = 37899.5
$ws.Range("J75").Value = 44437.668
$ws.Range("L75").Value = 44437.668
$ws.Range("N75").Value = -46309.668

# row 76
$ws.Range("H76").Value = 7169.091
$ws.Range("I76").Value = 6086.3
$ws.Range("K76").Value = 6086.3
$ws.Range("M76").Value = -5771.3

# row 78
$ws.Range("H78").Value = 37899.5
$ws.Range("J78").Value = 44437.668
$ws.Range("L78").Value = 133313.004
$ws.Range("N78").Value = -142673.004

# row 79
$ws.Range("H79").Value = 7169.091
$ws.Range("I79").Value = 6086.3
$ws.Range("K79").Value = 6086.3
$ws.Range("M79").Value = -4994.3

# row 98
$ws.Range("H98").Value = 2970.25
$ws.Range("I98").Value = 3058.2727
$ws.Range("K98").Value = 3058.2727
$ws.Range("M98").Value = -1560.2727

# row 100
$ws.Range("H100").Value = 3286.6
$ws.Range("I100").Value = 5233
$ws.Range("K100").Value = 5233
$ws.Range("M100").Value = -4692

# row 108
$ws.Range("H108").Value = 59343
$ws.Range("J108").Value = 59343
$ws.Range("L108").Value = 59343
$ws.Range("N108").Value = -67023

# row 109
$ws.Range("H109").Value = 46000
$ws.Range("J109").Value = 46000
$ws.Range("L109").Value = 46000
$ws.Range("N109").Value = -48774

# row 122
$ws.Range("H122").Value = 2970.25
$ws.Range("I122").Value = 3058.2727
$ws.Range("K122").Value = 9174.8181
$ws.Range("M122").Value = -6724.8181

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3131156.5
$ws.Range("I32").Value = 3283508.5
$ws.Range("K32").Value = 3283508.5
$ws.Range("M32").Value = -3283221.5

# row 74
$ws.Range("H74").Value = 19847.91
$ws.Range("I74").Value = 25672.658
$ws.Range("K74").Value = 25672.658
$ws.Range("M74").Value = -24798.658

# row 77
$ws.Range("H77").Value = 19847.91
$ws.Range("I77").Value = 25672.658
$ws.Range("K77").Value = 128363.29
$ws.Range("M77").Value = -123995.29

$ws = $wb.Worksheets.Item("BSM")
# row 75
$ws.Range("H75").Value = 15171
$ws.Range("I75").Value = 3756.5
$ws.Range("K75").Value = 3756.5
$ws.Range("M75").Value = -2820.5

# row 78
$ws.Range("H78").Value = 15171
$ws.Range("I78").Value = 3756.5
$ws.Range("K78").Value = 11269.5
$ws.Range("M78").Value = -6589.5

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 4442.933
$ws.Range("I16").Value = 646.6
$ws.Range("J16").Value = 8239.267
$ws.Range("K16").Value = 646.6
$ws.Range("L16").Value = 8239.267
$ws.Range("M16").Value = -359.6
$ws.Range("N16").Value = -8813.267

# row 31
$ws.Range("H31").Value = 4746.234
$ws.Range("J31").Value = 7334.2
$ws.Range("L31").Value = 7334.2
$ws.Range("N31").Value = -7924.2

# row 34
$ws.Range("H34").Value = 4746.234
$ws.Range("J34").Value = 7334.2
$ws.Range("L34").Value = 7334.2
$ws.Range("N34").Value = -7738.2

# row 113
$ws.Range("H113").Value = 4442.933
$ws.Range("I113").Value = 646.6
$ws.Range("J113").Value = 8239.267
$ws.Range("K113").Value = 646.6
$ws.Range("L113").Value = 8239.267
$ws.Range("M113").Value = 1523.4
$ws.Range("N113").Value = -12579.267

# row 122
$ws.Range("H122").Value = 4776.3335
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# row 134
$ws.Range("H134").Value = 4786.4346
$ws.Range("I134").Value = 2745.3
$ws.Range("J134").Value = 6356.5386
$ws.Range("K134").Value = 8235.900000000001
$ws.Range("L134").Value = 19069.6158
$ws.Range("M134").Value = -5700.900000000001
$ws.Range("N134").Value = -24139.6158

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1583.1428
$ws.Range("J5").Value = 4001.6667
$ws.Range("L5").Value = 12005.0001
$ws.Range("N5").Value = -12229.0001

# row 23
$ws.Range("H23").Value = 174.27272
$ws.Range("J23").Value = 272.16666
$ws.Range("L23").Value = 816.4999799999999
$ws.Range("N23").Value = -1286.49998

# row 98
$ws.Range("H98").Value = 1512.9
$ws.Range("J98").Value = 1722.125
$ws.Range("L98").Value = 5166.375
$ws.Range("N98").Value = -8162.375

# row 113
$ws.Range("H113").Value = 6021.923
$ws.Range("J113").Value = 6473.75
$ws.Range("L113").Value = 19421.25
$ws.Range("N113").Value = -23761.25

# row 135
$ws.Range("H135").Value = 1583.1428
$ws.Range("J135").Value = 4001.6667
$ws.Range("L135").Value = 36015.0003
$ws.Range("N135").Value = -41085.0003

$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1031.25
$ws.Range("I97").Value = 766.6667
$ws.Range("K97").Value = 766.6667
$ws.Range("M97").Value = -270.6667

# row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 9447.395500000001
$ws.Range("J136").Value = 17449.9
$ws.Range("L136").Value = 52349.7
$ws.Range("N136").Value = -57449.7

$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("J136").Value = 441276.56
$ws.Range("L136").Value = 1323829.68
$ws.Range("N136").Value = -1328929.68
